# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Thu Aug 31 13:21:12 UTC 2023 with GitHub Actions".
# D = Price, E = Volume(1h) columns, rows 2-51 (one row per coin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='27.342.73'},
    @{Cell='E2'; Value='  -0.60%  '},
    @{Cell='D3'; Value='1.715.78'},
    @{Cell='E3'; Value='  -0.53%  '},
    @{Cell='D4'; Value='1.005'},
    @{Cell='E4'; Value='  +0.25%  '},
    @{Cell='D5'; Value='224.61'},
    @{Cell='E5'; Value='  -0.31%  '},
    @{Cell='D6'; Value='0.5298'},
    @{Cell='E6'; Value='  -0.95%  '},
    @{Cell='D7'; Value='1.006'},
    @{Cell='E7'; Value='  +0.30%  '},
    @{Cell='D8'; Value='0.06705'},
    @{Cell='E8'; Value='  +1.51%  '},
    @{Cell='D9'; Value='0.2653'},
    @{Cell='E9'; Value='  -0.31%  '},
    @{Cell='D10'; Value='20.91'},
    @{Cell='E10'; Value='  -3.02%  '},
    @{Cell='D11'; Value='0.07682'},
    @{Cell='E11'; Value='  +0.23%  '},
    @{Cell='D12'; Value='4.491'},
    @{Cell='E12'; Value='  -2.34%  '},
    @{Cell='D13'; Value='1.952.61'},
    @{Cell='E13'; Value='  -0.47%  '},
    @{Cell='D14'; Value='1.722.26'},
    @{Cell='E14'; Value='  -0.13%  '},
    @{Cell='D15'; Value='0.5808'},
    @{Cell='E15'; Value='  +0.08%  '},
    @{Cell='D16'; Value='0.0₅8210'},
    @{Cell='E16'; Value='  -0.98%  '},
    @{Cell='D17'; Value='67.83'},
    @{Cell='E17'; Value='  +0.01%  '},
    @{Cell='D18'; Value='27.375.50'},
    @{Cell='E18'; Value='  -0.43%  '},
    @{Cell='D19'; Value='222.58'},
    @{Cell='E19'; Value='  +1.82%  '},
    @{Cell='D20'; Value='1.006'},
    @{Cell='E20'; Value='  +0.33%  '},
    @{Cell='D21'; Value='4.656'},
    @{Cell='E21'; Value='  -1.38%  '},
    @{Cell='D22'; Value='10.45'},
    @{Cell='E22'; Value='  -1.30%  '},
    @{Cell='D23'; Value='6.024'},
    @{Cell='E23'; Value='  -0.19%  '},
    @{Cell='E24'; Value='  +0.25%  '},
    @{Cell='D25'; Value='145.31'},
    @{Cell='E25'; Value='  +1.51%  '},
    @{Cell='D26'; Value='1.706'},
    @{Cell='E26'; Value='  -2.66%  '},
    @{Cell='D27'; Value='0.1207'},
    @{Cell='E27'; Value='  -2.09%  '},
    @{Cell='D28'; Value='7.250'},
    @{Cell='E28'; Value='  -1.10%  '},
    @{Cell='D29'; Value='16.24'},
    @{Cell='E29'; Value='  -1.47%  '},
    @{Cell='D30'; Value='0.05389'},
    @{Cell='E30'; Value='  -1.84%  '},
    @{Cell='D31'; Value='1.295'},
    @{Cell='E31'; Value='  -0.38%  '},
    @{Cell='D32'; Value='3.481'},
    @{Cell='E32'; Value='  -2.04%  '},
    @{Cell='D33'; Value='3.415'},
    @{Cell='E33'; Value='  -0.77%  '},
    @{Cell='D34'; Value='1.637'},
    @{Cell='E34'; Value='  -1.47%  '},
    @{Cell='E35'; Value='  +0.26%  '},
    @{Cell='D36'; Value='0.9527'},
    @{Cell='E36'; Value='  -0.50%  '},
    @{Cell='D37'; Value='2.393'},
    @{Cell='E37'; Value='  -1.28%  '},
    @{Cell='D38'; Value='0.5908'},
    @{Cell='E38'; Value='  -0.46%  '},
    @{Cell='D39'; Value='1.150.29'},
    @{Cell='E39'; Value='  +9.61%  '},
    @{Cell='D40'; Value='0.01654'},
    @{Cell='E40'; Value='  +0.34%  '},
    @{Cell='D41'; Value='5.842'},
    @{Cell='E41'; Value='  -0.94%  '},
    @{Cell='E42'; Value='  +0.31%  '},
    @{Cell='D43'; Value='0.8420'},
    @{Cell='E43'; Value='  -0.68%  '},
    @{Cell='D44'; Value='100.99'},
    @{Cell='E44'; Value='  -0.34%  '},
    @{Cell='D45'; Value='1.859.33'},
    @{Cell='E45'; Value='  -0.49%  '},
    @{Cell='E46'; Value='  +0.56%  '},
    @{Cell='D47'; Value='57.82'},
    @{Cell='E47'; Value='  -1.64%  '},
    @{Cell='D48'; Value='0.4583'},
    @{Cell='E48'; Value='  +2.28%  '},
    @{Cell='D49'; Value='8.135'},
    @{Cell='E49'; Value='  -0.93%  '},
    @{Cell='D50'; Value='1.002'},
    @{Cell='E50'; Value='  -0.03%  '},
    @{Cell='D51'; Value='0.05196'},
    @{Cell='E51'; Value='  -0.97%  '}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "1.005")
    # are not silently coerced into numbers by the Value setter.
    $r.NumberFormat = "@"
    $r.Value = $u.Value
    $r.Style = "Normal"
}
